$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the xpath selector values in column C (rows 4 and 5) to the new xpath
$newXPath = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[5]/a'
$ws.Range("C4").Value = $newXPath
$ws.Range("C5").Value = $newXPath

# Set column A width (approx. 22.42578125 characters as authored in Excel)
$ws.Columns.Item(1).ColumnWidth = 21.6

# Set the active selection to C3
$ws.Range("C3").Select()
